$wb = $excel.ActiveWorkbook

$sheetNames = @(
    @("血肉太白", "血肉太白"),
    @("血肉太白1", "血肉太白"),
    @("白煙", "白煙"),
    @("白煙1", "白煙"),
    @("白煙2", "白煙"),
    @("白煙3", "白煙"),
    @("白煙4", "白煙"),
    @("血肉太白2", "血肉太白"),
    @("白煙5", "白煙"),
    @("白煙6", "白煙"),
    @("白煙7", "白煙"),
    @("白煙8", "白煙"),
    @("白煙9", "白煙"),
    @("白煙10", "白煙"),
    @("白煙11", "白煙"),
    @("血肉太白3", "血肉太白"),
    @("白煙12", "白煙"),
    @("血肉太白4", "血肉太白"),
    @("白煙13", "白煙"),
    @("白煙14", "白煙"),
    @("血肉太白5", "血肉太白")
)

foreach ($pair in $sheetNames) {
    $tabName = $pair[0]
    $cellA1 = $pair[1]
    $count = $wb.Worksheets.Count
    $last = $wb.Worksheets.Item($count)
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
    $ws.Name = $tabName
    $ws.Range("A1").Value = $cellA1
    $ws.Range("B1").Value = "Start Frame #"
    $ws.Range("C1").Value = "Start Time"
    $ws.Range("D1").Value = "End Frame #"
    $ws.Range("E1").Value = "End Time"
}
